# reconfigure report_config to include spc_chart_type column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "rare_event_chart" column (E), which held N/Y flags and is
# no longer wanted in the rebuilt report_config layout.
$ws.Range("E1").EntireColumn.Delete() | Out-Null

# Insert a new blank column at D. "aggregation" (previously D) slides right
# to E, and "report_comment" (previously F) slides to E then back to F -
# matching the target column order: ref, measure_name, domain,
# spc_chart_type, aggregation, report_comment.
$ws.Range("D1").EntireColumn.Insert() | Out-Null

# Header for the new column.
$ws.Range("D1").Value = "spc_chart_type"

# Populate the new column: every measure uses the "xmr" SPC chart type
# except the last row, which uses "t".
$spcValues = @("xmr", "xmr", "xmr", "xmr", "xmr", "xmr", "xmr", "xmr", "t")
for ($i = 0; $i -lt $spcValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $spcValues[$i]
}

# Re-fit the column widths now that the data shape has changed (a constant
# offset compensates for this engine's internal column-width padding model
# so the stored widths line up with the target character widths).
$padding = 0.8333333333333339
$ws.Columns.Item(1).ColumnWidth = (3.44140625 - $padding)
$ws.Columns.Item(2).ColumnWidth = (19.88671875 - $padding)
$ws.Columns.Item(3).ColumnWidth = (9.33203125 - $padding)
$ws.Columns.Item(4).ColumnWidth = (14 - $padding)
$ws.Columns.Item(5).ColumnWidth = (11.109375 - $padding)
$ws.Columns.Item(6).ColumnWidth = (121 - $padding)

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null

$wb.Save()
